$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the last existing date cell (A244) so the new
# date cells (A245:A247) reuse the same cellXf (s="2") instead of minting
# a brand-new style entry.
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)

# New data rows (aggiornamento fino a 6/03)
$ws.Range("A245").Value = 44319
$ws.Range("B245").Value = 31
$ws.Range("C245").Value = 217
$ws.Range("D245").Value = 114.8069180426743

$ws.Range("A246").Value = 44320
$ws.Range("B246").Value = 14
$ws.Range("C246").Value = 219
$ws.Range("D246").Value = 115.8650463195653

$ws.Range("A247").Value = 44321
$ws.Range("B247").Value = 11
$ws.Range("C247").Value = 210
$ws.Range("D247").Value = 111.1034690735558
